$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Insert a new blank row above the current row 7 (AL population row),
# pushing "AL" down to row 8 and "AL % of USA" down to row 9.
$ws.Rows("7:7").Insert()
$ws.Range("A7").Clear()

# Add the new "Delta" (year-over-year AL population change) row at row 10.
$ws.Range("A10").Value = "Delta"
$ws.Range("A10").HorizontalAlignment = -4131

$ws.Range("C10").Formula = "=C8-B8"
$ws.Range("D10:K10").FormulaR1C1 = "=R[-2]C[0]-R[-2]C[-1]"
$ws.Range("C10:K10").NumberFormat = "#,##0"
$ws.Range("C10:K10").HorizontalAlignment = -4131

$ws.Range("E10").Select()
